# Append a new job listing row to the "ランサーズ" sheet and refresh the
# "取得日時" (fetched-at) timestamp on every existing row to the new run's
# timestamp, per commit message "Append: 2026-02-01 02:45 JST".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-02-01 02:45:04"

# Hyperlinks in this engine are not re-anchored by row Insert/shift
# operations, so drop them up front and rebuild the full set once the
# final row layout is in place (see below).
$ws.Hyperlinks.Delete()

# Push the existing data rows (old rows 3-5) down by one to make room for
# the new listing, which is inserted right after the top (most recent)
# existing row.
$ws.Rows.Item(3).Insert()

# New row 3: "Windows or Mac対応|本人顔ベースのリアルタイム顔変換システム開発"
$ws.Range("A3").Value = $newTimestamp
$ws.Range("B3").Value = "Windows or Mac対応|本人顔ベースのリアルタイム顔変換システム開発"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5483207"
$ws.Range("G3").Value = 125
$ws.Range("H3").Value = "◆開発,システム開発"

# Refresh the "取得日時" timestamp on every other row to the same run.
$ws.Range("A2").Value = $newTimestamp
$ws.Range("A4").Value = $newTimestamp
$ws.Range("A5").Value = $newTimestamp
$ws.Range("A6").Value = $newTimestamp

# Rebuild all hyperlinks on column F (rows 2-6) so the relationship ids
# line up with the new row order: row3 is the brand-new listing, rows
# 4-6 carry forward the URLs that used to live in rows 3-5.
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5482904")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5483207")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5482939")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5482932")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5482835")

# Hyperlinks.Add() re-styles the target cell with a fresh auto-generated
# xf; pin it back to the workbook's named "Hyperlink" cell style so the
# F-column look matches the rest of the sheet.
$ws.Range("F2").Style = "Hyperlink"
$ws.Range("F3").Style = "Hyperlink"
$ws.Range("F4").Style = "Hyperlink"
$ws.Range("F5").Style = "Hyperlink"
$ws.Range("F6").Style = "Hyperlink"

# Column B: 38 -> 41 ; Column D: 28 -> 30 (raw stored character widths).
# ColumnWidth read/write goes through Excel's padding-adjusted unit
# (offset by 5/6 of a character), so subtract that offset to land on the
# exact integer width that ends up serialized in the xlsx.
$ws.Columns.Item(2).ColumnWidth = 41 - (5/6)
$ws.Columns.Item(4).ColumnWidth = 30 - (5/6)
